# Auto-generated edit script: updates the cryptos price/volume table
# cells per the commit diff (prices, % changes, and the Dai/SuiNetwork
# row swap). Each touched cell is written as literal text (apostrophe-
# prefixed) so Excel's automatic number/date inference never mangles
# values such as "1.00", "69.075.04" or "0.0₃0821"; the quote-prefix
# marker left behind is then cleared by resetting the cell style back
# to "Normal" so no stray style attribute is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = ('''69.259.61') }
    @{ Cell = "E2"; Value = ('''  -0.33%  ') }
    @{ Cell = "D3"; Value = ('''2.471.41') }
    @{ Cell = "E3"; Value = ('''  -0.85%  ') }
    @{ Cell = "E4"; Value = ('''  -0.07%  ') }
    @{ Cell = "D5"; Value = ('''559.92') }
    @{ Cell = "E5"; Value = ('''  -1.57%  ') }
    @{ Cell = "D6"; Value = ('''163.67') }
    @{ Cell = "E6"; Value = ('''  -1.71%  ') }
    @{ Cell = "E7"; Value = ('''  +0.00%  ') }
    @{ Cell = "E8"; Value = ('''  -0.86%  ') }
    @{ Cell = "D9"; Value = ('''2.469.02') }
    @{ Cell = "E9"; Value = ('''  -0.91%  ') }
    @{ Cell = "D10"; Value = ('''0.153') }
    @{ Cell = "E10"; Value = ('''  -3.61%  ') }
    @{ Cell = "E11"; Value = ('''  -0.56%  ') }
    @{ Cell = "D12"; Value = ('''0.336') }
    @{ Cell = "E12"; Value = ('''  -4.31%  ') }
    @{ Cell = "D13"; Value = ('''4.84') }
    @{ Cell = "E13"; Value = ('''  -0.74%  ') }
    @{ Cell = "D14"; Value = ('''2.922.89') }
    @{ Cell = "E14"; Value = ('''  -0.94%  ') }
    @{ Cell = "D15"; Value = ('''68.880.93') }
    @{ Cell = "E15"; Value = ('''  -0.65%  ') }
    @{ Cell = "E16"; Value = ('''  -2.63%  ') }
    @{ Cell = "D17"; Value = ('''23.67') }
    @{ Cell = "E17"; Value = ('''  -2.26%  ') }
    @{ Cell = "D18"; Value = ('''2.484.02') }
    @{ Cell = "E18"; Value = ('''  -1.35%  ') }
    @{ Cell = "D19"; Value = ('''10.79') }
    @{ Cell = "E19"; Value = ('''  -3.79%  ') }
    @{ Cell = "D20"; Value = ('''343.29') }
    @{ Cell = "E20"; Value = ('''  -2.55%  ') }
    @{ Cell = "D21"; Value = ('''7.08') }
    @{ Cell = "E21"; Value = ('''  -4.00%  ') }
    @{ Cell = "D22"; Value = ('''3.82') }
    @{ Cell = "E22"; Value = ('''  -2.19%  ') }
    @{ Cell = "E23"; Value = ('''  -0.52%  ') }
    @{ Cell = "B24"; Value = ('''SuiNetwork') }
    @{ Cell = "C24"; Value = ('''https://coinranking.com/coin/3xJluUMvp+suinetwork-sui') }
    @{ Cell = "D24"; Value = ('''1.91') }
    @{ Cell = "E24"; Value = ('''  +0.58%  ') }
    @{ Cell = "B25"; Value = ('''Dai') }
    @{ Cell = "C25"; Value = ('''https://coinranking.com/coin/MoTuySvg7+dai-dai') }
    @{ Cell = "D25"; Value = ('''1.00') }
    @{ Cell = "E25"; Value = ('''  +0.03%  ') }
    @{ Cell = "D26"; Value = ('''67.24') }
    @{ Cell = "E26"; Value = ('''  -3.05%  ') }
    @{ Cell = "D27"; Value = ('''3.71') }
    @{ Cell = "E27"; Value = ('''  -2.36%  ') }
    @{ Cell = "E28"; Value = ('''  -0.95%  ') }
    @{ Cell = "D29"; Value = ('''0.997') }
    @{ Cell = "E29"; Value = ('''  -0.23%  ') }
    @{ Cell = "D30"; Value = ('''8.21') }
    @{ Cell = "E30"; Value = ('''  -4.82%  ') }
    @{ Cell = "D31"; Value = ('''0.0' + [char]0x2083 + '0823') }
    @{ Cell = "E31"; Value = ('''  -5.58%  ') }
    @{ Cell = "D32"; Value = ('''7.21') }
    @{ Cell = "E32"; Value = ('''  -5.08%  ') }
    @{ Cell = "D33"; Value = ('''442.08') }
    @{ Cell = "E33"; Value = ('''  -0.25%  ') }
    @{ Cell = "E34"; Value = ('''  -0.06%  ') }
    @{ Cell = "E35"; Value = ('''  -3.48%  ') }
    @{ Cell = "D36"; Value = ('''1.62') }
    @{ Cell = "E36"; Value = ('''  -5.13%  ') }
    @{ Cell = "D37"; Value = ('''157.41') }
    @{ Cell = "E37"; Value = ('''  +2.03%  ') }
    @{ Cell = "D38"; Value = ('''19.07') }
    @{ Cell = "E38"; Value = ('''  +0.01%  ') }
    @{ Cell = "E39"; Value = ('''  -0.01%  ') }
    @{ Cell = "E40"; Value = ('''  -3.10%  ') }
    @{ Cell = "D41"; Value = ('''17.95') }
    @{ Cell = "E41"; Value = ('''  -1.02%  ') }
    @{ Cell = "E42"; Value = ('''  -3.06%  ') }
    @{ Cell = "D43"; Value = ('''4.48') }
    @{ Cell = "E43"; Value = ('''  -2.52%  ') }
    @{ Cell = "D44"; Value = ('''37.45') }
    @{ Cell = "E44"; Value = ('''  -0.89%  ') }
    @{ Cell = "D45"; Value = ('''1.49') }
    @{ Cell = "E45"; Value = ('''  -5.94%  ') }
    @{ Cell = "E46"; Value = ('''  +3.25%  ') }
    @{ Cell = "D47"; Value = ('''2.09') }
    @{ Cell = "E47"; Value = ('''  -4.15%  ') }
    @{ Cell = "D48"; Value = ('''133.73') }
    @{ Cell = "E48"; Value = ('''  -3.53%  ') }
    @{ Cell = "D49"; Value = ('''3.37') }
    @{ Cell = "E49"; Value = ('''  -1.85%  ') }
    @{ Cell = "D50"; Value = ('''0.0719') }
    @{ Cell = "E50"; Value = ('''  -0.53%  ') }
    @{ Cell = "D51"; Value = ('''0.486') }
    @{ Cell = "E51"; Value = ('''  -3.71%  ') }
)

foreach ($u in $updates) {
    $c = $ws.Range($u.Cell)
    $c.Value = $u.Value
    $c.Style = "Normal"
}

